$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price reading was inserted as row 305 (Primera quality, Cilantro,
# Terminal Hortofrutícola Agro Chillán), pushing the existing rows 305-357 down
# to 306-358.
$ws.Rows.Item(305).Insert()

$ws.Cells.Item(305, 1).Value = 7
$ws.Cells.Item(305, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(305, 3).Value = "Ñuble"
$ws.Cells.Item(305, 4).Value = 45258
$ws.Cells.Item(305, 5).Value = 16
$ws.Cells.Item(305, 6).Value = 100112040
$ws.Cells.Item(305, 7).Value = "Cilantro"
$ws.Cells.Item(305, 8).Value = "Sin especificar"
$ws.Cells.Item(305, 9).Value = "Primera"
$ws.Cells.Item(305, 10).Value = 150
$ws.Cells.Item(305, 11).Value = 2000
$ws.Cells.Item(305, 12).Value = 2000
$ws.Cells.Item(305, 13).Value = 2000
$ws.Cells.Item(305, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(305, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(305, 16).Value = 2000
$ws.Cells.Item(305, 17).Value = 1
$ws.Cells.Item(305, 18).Value = "Hortaliza"

# Match the date-formatted style used by the rest of column D.
$ws.Cells.Item(305, 4).NumberFormat = $ws.Cells.Item(306, 4).NumberFormat
